# This script reproduces the commit that re-ran the differential-privacy
# logistic-regression simulation ("numba is ready") and, as a side effect of
# how the shared-string table was rebuilt, moved the "Weights" column header
# from column C to column K (last) on both worksheets, with every other header
# shifting one column to the left. All numeric results were also recomputed.

$wb = $excel.ActiveWorkbook

$headers = @('N', '$\epsilon = 0.0005$', '$\epsilon = 0.001$', '$\epsilon = 0.005$', '$\epsilon = 0.01$', '$\epsilon = 0.05$', '$\epsilon = 0.1$', '$\epsilon = 0.5$', '$\epsilon = 3$', 'Weights')
$cols = @('B', 'C', 'D', 'E', 'F', 'G', 'H', 'I', 'J', 'K')

# ---- Worksheet 1 ----
$ws = $wb.Worksheets.Item(1)

# Rewrite the header row so the columns read N, epsilon=0.0005 .. epsilon=3, Weights
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "1").Value() = $headers[$i]
}

# Recomputed simulation output: row layout is N, then the 9 data columns
# (C..K) in the new column order (epsilon=0.0005 .. epsilon=3, Weights)
$rows = @(
    @(2, @('21', '158.0688278697604', '79.05336545289872', '15.77442546854861', '7.911110477161534', '1.579329682280626', '0.7893463099091917', '0.1578901368217719', '0.02629573490310636', '0.001637830735344022')),
    @(3, @('42', '157.9936466259097', '78.86714290354189', '15.84026788746038', '7.896287163348887', '1.578787305226836', '0.7904772314427679', '0.1577971005734869', '0.02635086897430705', '0.001263670303085209')),
    @(4, @('63', '158.0106332225969', '79.1854257510394', '15.8079945778298', '7.909653176262839', '1.575051601553223', '0.788970150746177', '0.1579068543628789', '0.02643033426596556', '0.00125661083996006')),
    @(5, @('84', '157.8264703766389', '79.11088157410612', '15.75467045057052', '7.899593005543904', '1.580623648535944', '0.7924772333903534', '0.1583236368438386', '0.02631171628035118', '0.001172953421322764')),
    @(6, @('105', '158.3703944047845', '78.95949697530881', '15.81737603425909', '7.907589818270206', '1.581647538379609', '0.7914178172239767', '0.1581269137705753', '0.02634953162634622', '0.0004758824762643964')),
    @(7, @('126', '157.819168518465', '79.03703695391792', '15.81620235133372', '7.907834739957322', '1.581445677609782', '0.7898172878621201', '0.1581725699100923', '0.02636556072963221', '0.001283062887174413')),
    @(8, @('147', '157.7064334926939', '79.14275146460395', '15.79859531894928', '7.90817213552601', '1.578255476177183', '0.7897679477538857', '0.1577045552998541', '0.02634600115552062', '0.0002468555559240835')),
    @(9, @('168', '158.077046910276', '79.00182780034851', '15.7893807605487', '7.908362050662186', '1.581679992131252', '0.7887079040819422', '0.1582168938604664', '0.02631528379336331', '8.433435698692454e-05')),
    @(10, @('189', '158.1648500454752', '78.96925490385999', '15.81272143428852', '7.914238040756166', '1.581450900658315', '0.790035367734804', '0.1582959028008062', '0.02630744535160921', '-0.000347389520512437')),
    @(11, @('210', '158.2386448695727', '78.80097424237991', '15.80663277067362', '7.908449611390294', '1.581340383426711', '0.7899110040838146', '0.1578444353889548', '0.02630472123321671', '-0.0001551035771935032')),
    @(12, @('231', '158.0621574291321', '79.16050449892978', '15.79237389519561', '7.917799444688045', '1.579782560716009', '0.7906671902218425', '0.1583295872741402', '0.02636122437006319', '-0.0001906093313675458')),
    @(13, @('252', '158.0987054313068', '78.90388884650565', '15.77671768711285', '7.925496359190459', '1.577241241926728', '0.7899001872789395', '0.1578590172426838', '0.02629645447558512', '-0.0005828888420591587')),
    @(14, @('273', '157.9517182604025', '78.84022008118109', '15.77381847377514', '7.900640635118423', '1.580605594325018', '0.7917953462130278', '0.1586813218178441', '0.02636423608493736', '0.00087430426807908')),
    @(15, @('294', '157.7258953840803', '78.88333445460057', '15.84159065333031', '7.907640730730147', '1.579106709261391', '0.7896644720224095', '0.1580550011812424', '0.02634704698199194', '0.001627840560397385')),
    @(16, @('315', '158.089302378476', '79.06467286625221', '15.81367262652027', '7.911791930586967', '1.582015793768093', '0.7903568898610356', '0.1579501435514131', '0.02642692871810127', '0.00161531834432339'))
)

foreach ($entry in $rows) {
    $r = $entry[0]
    $vals = $entry[1]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = $cols[$i] + $r
        if ($i -eq 0) {
            $ws.Range($addr).Value() = [int]$vals[$i]
        } else {
            $ws.Range($addr).Value() = [double]$vals[$i]
        }
    }
}

# ---- Worksheet 2 ----
$ws = $wb.Worksheets.Item(2)

# Rewrite the header row so the columns read N, epsilon=0.0005 .. epsilon=3, Weights
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "1").Value() = $headers[$i]
}

# Recomputed simulation output: row layout is N, then the 9 data columns
# (C..K) in the new column order (epsilon=0.0005 .. epsilon=3, Weights)
$rows = @(
    @(2, @('21', '0.5388971139192631', '0.1664125465440106', '0.002182387485425093', '0.001503397383743446', '1.634487692138122e-05', '1.011565171622803e-05', '5.882899279185047e-07', '6.396663408410638e-09', '0.003332138915318018')),
    @(3, @('42', '0.3335323798321853', '0.03572147323037055', '0.002475247567155695', '0.0004035365706866905', '2.643348702667772e-05', '6.427989801974283e-06', '4.45252362129621e-07', '4.214114203172238e-09', '0.003767373189123533')),
    @(4, @('63', '0.3008125773376262', '0.03514962242569448', '0.001876222959648762', '0.0005838257044296208', '3.079697062622574e-05', '1.171617735211757e-05', '2.910251648569854e-07', '4.471074262604271e-09', '0.004404665111276523')),
    @(5, @('84', '0.01932884489126172', '0.1006075961147395', '0.0038592891592249', '0.001455113291182925', '1.554817879565704e-05', '1.269924375156924e-05', '4.247439195904293e-07', '1.211004136715657e-08', '0.004701543263375125')),
    @(6, @('105', '0.1575168639765247', '0.240170097911845', '0.004476703788686597', '0.001206297079026628', '9.719829027587821e-06', '1.170915032022572e-05', '2.424456167618441e-07', '8.950391493102031e-09', '0.005311560451219446')),
    @(7, @('126', '0.1337965642494017', '0.07880443229547093', '0.002355992711585806', '0.0008650009184373336', '4.139375405241967e-05', '4.856204543240151e-06', '2.428254368799974e-07', '6.895632761711884e-09', '0.001258530223895063')),
    @(8, @('147', '0.5160521199099239', '0.1050571029393607', '0.00390920208018836', '0.000635778347066127', '4.883158279649392e-06', '6.357338445237126e-06', '2.364545991994919e-07', '7.323844590964111e-09', '0.005670030180593685')),
    @(9, @('168', '0.08004406598556439', '0.1130714190438923', '0.00509915606076483', '0.001665191674019315', '1.692227814074824e-05', '1.006790094344584e-05', '3.121669004591994e-07', '1.055792321325364e-08', '0.005781157698736042')),
    @(10, @('189', '0.710564794046785', '0.07238609733290943', '0.004281897183314638', '0.0007026731990371843', '5.554532040487149e-05', '6.947884146959398e-06', '5.999492166946845e-07', '9.577691802518973e-09', '0.006049668660252355')),
    @(11, @('210', '0.07546626211931667', '0.03105487268236645', '0.002282794550209682', '0.001335590122854807', '2.226516105539521e-05', '6.254092754772214e-06', '6.222764334533833e-07', '1.929645180875972e-09', '0.006300797354183731')),
    @(12, @('231', '0.3259830440040595', '0.09396517814153946', '0.003993786168457676', '0.002559544451310559', '1.919523306323559e-05', '2.740140281862689e-05', '1.488700388683028e-07', '1.514753493006979e-08', '0.006499712033847863')),
    @(13, @('252', '0.2246792436929469', '0.0939191814506338', '0.001355182031456266', '0.0003458398931492565', '2.749801794867621e-05', '1.231391417578391e-05', '1.920299593956377e-07', '1.328687565981286e-08', '0.006870055417105785')),
    @(14, @('273', '0.3056032949815163', '0.03714426831907866', '0.006883794923718129', '0.0006258660675599264', '7.818105843218028e-05', '9.675663969027826e-06', '3.119007846344517e-07', '1.236937274783853e-08', '0.001984140473647917')),
    @(15, @('294', '0.256661068438039', '0.1344784169910012', '0.00458523391110811', '0.0003324442987628269', '3.676226086865599e-05', '1.019010566417569e-05', '1.23503459029077e-07', '1.674979616051067e-08', '0.002403459947195086')),
    @(16, @('315', '0.1142274935728161', '0.05448332624284676', '0.0029566651686651', '0.0003877569262783366', '2.717115401400796e-05', '1.33918030892793e-05', '5.822055549523775e-07', '4.537186840774366e-09', '0.002859514089584859'))
)

foreach ($entry in $rows) {
    $r = $entry[0]
    $vals = $entry[1]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = $cols[$i] + $r
        if ($i -eq 0) {
            $ws.Range($addr).Value() = [int]$vals[$i]
        } else {
            $ws.Range($addr).Value() = [double]$vals[$i]
        }
    }
}

